$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 336 - this shifts the existing rows 336:381 down to 337:382
# (matches the dimension change from A1:T381 to A1:T382 in the diff)
$ws.Rows(336).Insert()

# Populate the newly inserted row 336 with the new record's data.
$ws.Cells.Item(336, 1).Value = 7
$ws.Cells.Item(336, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(336, 3).Value = "Ñuble"
$ws.Cells.Item(336, 4).Value = 45212
$ws.Cells.Item(336, 5).Value = 16
$ws.Cells.Item(336, 6).Value = "Fruta"
$ws.Cells.Item(336, 7).Value = 100108
$ws.Cells.Item(336, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(336, 9).Value = 100108005
$ws.Cells.Item(336, 10).Value = "Piña"
$ws.Cells.Item(336, 11).Value = "Caramelo"
$ws.Cells.Item(336, 12).Value = "Segunda"
$ws.Cells.Item(336, 13).Value = 100
$ws.Cells.Item(336, 14).Value = 22000
$ws.Cells.Item(336, 15).Value = 23000
$ws.Cells.Item(336, 16).Value = 22500
$ws.Cells.Item(336, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(336, 18).Value = "Ecuador"
$ws.Cells.Item(336, 19).Value = 1607
$ws.Cells.Item(336, 20).Value = 14
